$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
Write-Host $ws.Name
